$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the active cell on the column we're about to remove (matches the
# selection left behind by the author after deleting the column), then
# delete the entire column M. This shifts column N (and everything to its
# right) one position to the left, so the former N-column values become
# the new column M, and the sheet's used range shrinks from N119 to M119.
$ws.Range("M1").Select()
$ws.Range("M:M").Delete()
